# Update EIA Table 4.5 (Industrial Sector) from "October 2016" vintage to
# "November 2016" vintage: refresh the title / rolling-12-months caption,
# insert the new "November" monthly data row, and refresh the Year-to-Date
# and Annual-Totals figures for 2014/2015/2016.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Title text and "Rolling 12 Months Ending in ..." caption
# ---------------------------------------------------------------------
$ws.Cells.Item(1,1).Value = "Table 4.5. Receipts, Average Cost, and Quality of Fossil Fuels: Industrial Sector, 2006 - November 2016 (continued)"
$ws.Cells.Item(57,1).Value = "Rolling 12 Months Ending in November"

# ---------------------------------------------------------------------
# 2) Insert the new "November" data row just above the "Year to Date"
#    section header (old row 53), push everything below it down by one.
#    Grab the number formatting from the row above (the December data
#    row, row 52) so the new row carries the normal data-row styles
#    instead of Excel's auto-guessed insert formatting.
# ---------------------------------------------------------------------
$ws.Rows.Item(53).Insert()
$ws.Range("A52:M52").Copy()
$ws.Range("A53:M53").PasteSpecial(-4122)

$ws.Cells.Item(53,1).Value = "November"
$ws.Cells.Item(53,2).Value = 200
$ws.Cells.Item(53,3).Value = 8
$ws.Cells.Item(53,4).Value = "W"
$ws.Cells.Item(53,5).Value = "W"
$ws.Cells.Item(53,6).Value = 5.47
$ws.Cells.Item(53,7).Value = 9
$ws.Cells.Item(53,8).Value = 65021
$ws.Cells.Item(53,9).Value = 63167
$ws.Cells.Item(53,10).Value = "W"
$ws.Cells.Item(53,11).Value = "W"
$ws.Cells.Item(53,12).Value = 59.8
$ws.Cells.Item(53,13).Value = "W"

# ---------------------------------------------------------------------
# 3) "Year to Date" block (now rows 55-57: 2014, 2015, 2016) - refresh
#    with the updated values.
# ---------------------------------------------------------------------
# 2014
$ws.Cells.Item(55,2).Value = 8450
$ws.Cells.Item(55,3).Value = 311
$ws.Cells.Item(55,6).Value = 5.82
$ws.Cells.Item(55,7).Value = 22
$ws.Cells.Item(55,8).Value = 677899
$ws.Cells.Item(55,9).Value = 656198
$ws.Cells.Item(55,12).Value = 62.8

# 2015
$ws.Cells.Item(56,1).Value = 2015
$ws.Cells.Item(56,2).Value = 7385
$ws.Cells.Item(56,3).Value = 274
$ws.Cells.Item(56,4).Value = "W"
$ws.Cells.Item(56,5).Value = "W"
$ws.Cells.Item(56,6).Value = 5.5
$ws.Cells.Item(56,7).Value = 23.4
$ws.Cells.Item(56,8).Value = 696317
$ws.Cells.Item(56,9).Value = 673448
$ws.Cells.Item(56,10).Value = "W"
$ws.Cells.Item(56,11).Value = "W"
$ws.Cells.Item(56,12).Value = 60.6
$ws.Cells.Item(56,13).Value = "W"

# 2016
$ws.Cells.Item(57,1).Value = 2016
$ws.Cells.Item(57,2).Value = 3031
$ws.Cells.Item(57,3).Value = 112
$ws.Cells.Item(57,4).Value = "W"
$ws.Cells.Item(57,5).Value = "W"
$ws.Cells.Item(57,6).Value = 5.83
$ws.Cells.Item(57,7).Value = 10.7
$ws.Cells.Item(57,8).Value = 684423
$ws.Cells.Item(57,9).Value = 663461
$ws.Cells.Item(57,10).Value = "W"
$ws.Cells.Item(57,11).Value = "W"
$ws.Cells.Item(57,12).Value = 57.8
$ws.Cells.Item(57,13).Value = "W"

# ---------------------------------------------------------------------
# 4) "Annual Totals / Rolling 12 Months" block (now rows 59-60: 2015,
#    2016) - refresh with the updated values.
# ---------------------------------------------------------------------
# 2015
$ws.Cells.Item(59,2).Value = 8670
$ws.Cells.Item(59,3).Value = 321
$ws.Cells.Item(59,6).Value = 5.55
$ws.Cells.Item(59,7).Value = 24.7
$ws.Cells.Item(59,8).Value = 760766
$ws.Cells.Item(59,9).Value = 735609
$ws.Cells.Item(59,12).Value = 60.7

# 2016
$ws.Cells.Item(60,2).Value = 3835
$ws.Cells.Item(60,3).Value = 142
$ws.Cells.Item(60,6).Value = 5.75
$ws.Cells.Item(60,7).Value = 12.5
$ws.Cells.Item(60,8).Value = 754070
$ws.Cells.Item(60,9).Value = 730989
$ws.Cells.Item(60,12).Value = 58.1
